$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = '(''reddit'', ''kotakuinact'', ''comment'')'
$ws.Range("C2").Value = 13
$ws.Range("B3").Value = '(''sea'', ''level'', ''rise'')'
$ws.Range("C3").Value = 7
$ws.Range("B4").Value = '(''climat'', ''chang'', ''polit'')'
$ws.Range("C4").Value = 6
$ws.Range("B5").Value = '(''giss'', ''nasa'', ''gov'')'
$ws.Range("C5").Value = 5
$ws.Range("B6").Value = '(''climat'', ''chang'', ''real'')'
$ws.Range("C6").Value = 5
$ws.Range("B7").Value = '(''unit'', ''nation'', ''publish'')'
$ws.Range("C7").Value = 4
$ws.Range("B8").Value = '(''nation'', ''publish'', ''report'')'
$ws.Range("C8").Value = 4
$ws.Range("B9").Value = '(''publish'', ''report'', ''last'')'
$ws.Range("C9").Value = 4
$ws.Range("B10").Value = '(''report'', ''last'', ''month'')'
$ws.Range("C10").Value = 4
$ws.Range("B11").Value = '(''last'', ''month'', ''warn'')'
$ws.Range("C11").Value = 4
$ws.Range("B12").Value = '(''month'', ''warn'', ''drastic'')'
$ws.Range("C12").Value = 4
$ws.Range("B13").Value = '(''warn'', ''drastic'', ''refor'')'
$ws.Range("C13").Value = 4
$ws.Range("B14").Value = '(''drastic'', ''refor'', ''effort'')'
$ws.Range("C14").Value = 4
$ws.Range("B15").Value = '(''refor'', ''effort'', ''reduct'')'
$ws.Range("C15").Value = 4
$ws.Range("B16").Value = '(''effort'', ''reduct'', ''carbon'')'
$ws.Range("C16").Value = 4
$ws.Range("B17").Value = '(''reduct'', ''carbon'', ''pollut'')'
$ws.Range("C17").Value = 4
$ws.Range("B18").Value = '(''carbon'', ''pollut'', ''planet'')'
$ws.Range("C18").Value = 4
$ws.Range("B19").Value = '(''pollut'', ''planet'', ''soon'')'
$ws.Range("C19").Value = 4
$ws.Range("B20").Value = '(''planet'', ''soon'', ''face'')'
$ws.Range("C20").Value = 4
$ws.Range("B21").Value = '(''soon'', ''face'', ''irrever'')'
$ws.Range("C21").Value = 4
$ws.Range("B22").Value = '(''face'', ''irrever'', ''global'')'
$ws.Range("C22").Value = 4
$ws.Range("B23").Value = '(''irrever'', ''global'', ''catastroph'')'
$ws.Range("C23").Value = 4
$ws.Range("B24").Value = '(''web'', ''archiv'', ''web'')'
$ws.Range("C24").Value = 4
$ws.Range("B25").Value = '(''climat'', ''chang'', ''go'')'
$ws.Range("C25").Value = 4
$ws.Range("B26").Value = '(''nasa'', ''gov'', ''tmp'')'
$ws.Range("C26").Value = 4
$ws.Range("B27").Value = '(''gov'', ''tmp'', ''gistemp'')'
$ws.Range("C27").Value = 4
$ws.Range("B28").Value = '(''tmp'', ''gistemp'', ''custom'')'
$ws.Range("C28").Value = 4
$ws.Range("B29").Value = '(''gistemp'', ''custom'', ''graph'')'
$ws.Range("C29").Value = 4
$ws.Range("B30").Value = '(''peopl'', ''3rd'', ''world'')'
$ws.Range("C30").Value = 4
$ws.Range("B31").Value = '(''unit'', ''state'', ''recov'')'
$ws.Range("C31").Value = 3
$ws.Range("B32").Value = '(''histor'', ''wildfir'', ''california'')'
$ws.Range("C32").Value = 3
$ws.Range("B33").Value = '(''esrl'', ''noaa'', ''gov'')'
$ws.Range("C33").Value = 3
$ws.Range("B34").Value = '(''mbrol67'', ''reddit'', ''kotakuinact'')'
$ws.Range("C34").Value = 3
$ws.Range("B35").Value = '(''scienc'', ''sciencemag'', ''content'')'
$ws.Range("C35").Value = 3
$ws.Range("B36").Value = '(''kotakuinact'', ''comment'', ''ethicssocju'')'
$ws.Range("C36").Value = 3
$ws.Range("B37").Value = '(''comment'', ''ethicssocju'', ''eurogam'')'
$ws.Range("C37").Value = 3
$ws.Range("B38").Value = '(''ethicssocju'', ''eurogam'', ''civil'')'
$ws.Range("C38").Value = 3
$ws.Range("B39").Value = '(''eurogam'', ''civil'', ''gather'')'
$ws.Range("C39").Value = 3
$ws.Range("B40").Value = '(''civil'', ''gather'', ''ea7hvkw'')'
$ws.Range("C40").Value = 3
$ws.Range("B41").Value = '(''climat'', ''chang'', ''alway'')'
$ws.Range("C41").Value = 3
$ws.Range("B42").Value = '(''due'', ''climat'', ''chang'')'
$ws.Range("C42").Value = 3
$ws.Range("B43").Value = '(''chang'', ''polit'', ''statement'')'
$ws.Range("C43").Value = 3
$ws.Range("B44").Value = '(''believ'', ''climat'', ''chang'')'
$ws.Range("C44").Value = 3
$ws.Range("B45").Value = '(''climat'', ''chang'', ''happen'')'
$ws.Range("C45").Value = 3
$ws.Range("B46").Value = '(''climat'', ''chang'', ''actual'')'
$ws.Range("C46").Value = 3
$ws.Range("B47").Value = '(''nuclear'', ''power'', ''plant'')'
$ws.Range("C47").Value = 3
$ws.Range("B48").Value = '(''believ'', ''man'', ''make'')'
$ws.Range("C48").Value = 3
$ws.Range("B49").Value = '(''graph'', ''png'', ''data'')'
$ws.Range("C49").Value = 3
$ws.Range("B50").Value = '(''png'', ''data'', ''giss'')'
$ws.Range("C50").Value = 3
$ws.Range("B51").Value = '(''data'', ''giss'', ''nasa'')'
$ws.Range("C51").Value = 3
$ws.Range("B52").Value = '(''climat'', ''chang'', ''scienc'')'
$ws.Range("C52").Value = 3
$ws.Range("B53").Value = '(''ice'', ''cap'', ''melt'')'
$ws.Range("C53").Value = 3
$ws.Range("B54").Value = '(''nation'', ''secur'', ''threat'')'
$ws.Range("C54").Value = 3
$ws.Range("B55").Value = '(''vike'', ''grow'', ''wheat'')'
$ws.Range("C55").Value = 3
$ws.Range("B56").Value = '(''absorb'', ''low'', ''frequenc'')'
$ws.Range("C56").Value = 3
$ws.Range("B57").Value = '(''3rd'', ''world'', ''west'')'
$ws.Range("C57").Value = 3
$ws.Range("B58").Value = '(''import'', ''peopl'', ''3rd'')'
$ws.Range("C58").Value = 3
$ws.Range("B59").Value = '(''tell'', ''3rd'', ''world'')'
$ws.Range("C59").Value = 3
$ws.Range("B60").Value = '(''3rd'', ''world'', ''breed'')'
$ws.Range("C60").Value = 3
$ws.Range("B61").Value = '(''world'', ''breed'', ''rabbit'')'
$ws.Range("C61").Value = 3
$ws.Range("B62").Value = '(''low'', ''iq'', ''peopl'')'
$ws.Range("C62").Value = 3
$ws.Range("B63").Value = '(''economi'', ''climat'', ''chang'')'
$ws.Range("C63").Value = 3
$ws.Range("B64").Value = '(''dislik'', ''big'', ''govern'')'
$ws.Range("C64").Value = 3
$ws.Range("B65").Value = '(''use'', ''coal'', ''oil'')'
$ws.Range("C65").Value = 3
$ws.Range("B66").Value = '(''eurogam'', ''shoehorn'', ''politic'')'
$ws.Range("C66").Value = 2
$ws.Range("B67").Value = '(''mile'', ''away'', ''sit'')'
$ws.Range("C67").Value = 2
$ws.Range("B68").Value = '(''away'', ''sit'', ''write'')'
$ws.Range("C68").Value = 2
$ws.Range("B69").Value = '(''sit'', ''write'', ''compani'')'
$ws.Range("C69").Value = 2
$ws.Range("B70").Value = '(''write'', ''compani'', ''restart'')'
$ws.Range("C70").Value = 2
$ws.Range("B71").Value = '(''compani'', ''restart'', ''gas'')'
$ws.Range("C71").Value = 2
$ws.Range("B72").Value = '(''restart'', ''gas'', ''frack'')'
$ws.Range("C72").Value = 2
$ws.Range("B73").Value = '(''first'', ''time'', ''sinc'')'
$ws.Range("C73").Value = 2
$ws.Range("B74").Value = '(''time'', ''sinc'', ''ban'')'
$ws.Range("C74").Value = 2
$ws.Range("B75").Value = '(''sinc'', ''ban'', ''unit'')'
$ws.Range("C75").Value = 2
$ws.Range("B76").Value = '(''ban'', ''unit'', ''kingdom'')'
$ws.Range("C76").Value = 2
$ws.Range("B77").Value = '(''unit'', ''kingdom'', ''incred'')'
$ws.Range("C77").Value = 2
$ws.Range("B78").Value = '(''kingdom'', ''incred'', ''frack'')'
$ws.Range("C78").Value = 2
$ws.Range("B79").Value = '(''incred'', ''frack'', ''begin'')'
$ws.Range("C79").Value = 2
$ws.Range("B80").Value = '(''frack'', ''begin'', ''exact'')'
$ws.Range("C80").Value = 2
$ws.Range("B81").Value = '(''begin'', ''exact'', ''week'')'
$ws.Range("C81").Value = 2
